$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3968591891046255
$ws.Range("C2").Value = 0.3932745487300524

$ws.Range("B3").Value = 46.08271915157921
$ws.Range("C3").Value = 46.46138267323168

$ws.Range("B4").Value = 773.0656828242987
$ws.Range("C4").Value = 769.8464956551795

$ws.Range("B5").Value = 72.76409810354545
$ws.Range("C5").Value = 68.55733984940724

$ws.Range("B6").Value = 25935.8031076495
$ws.Range("C6").Value = 23758.85461543648

$ws.Range("B7").Value = 953.4433762394547
$ws.Range("C7").Value = 1313.62912996436

$ws.Range("B8").Value = -2333.170494236158
$ws.Range("C8").Value = 98.07242789251477

$ws.Range("B9").Value = 582.295169670233
$ws.Range("C9").Value = 567.2022375458338

$ws.Range("B10").Value = 3360.526318588844
$ws.Range("C10").Value = 4179.677015473537

$ws.Range("B11").Value = -1455.95198368039
$ws.Range("C11").Value = -403.9338758972591

$ws.Range("B12").Value = -8.438350680131276
$ws.Range("C12").Value = -8.406545210008485

$ws.Range("B13").Value = -4.345022509496995
$ws.Range("C13").Value = -4.493072200383833

$ws.Range("B14").Value = -4.283884224258766
$ws.Range("C14").Value = -4.433297209417767

$ws.Range("B15").Value = -0.9031443591398607
$ws.Range("C15").Value = -1.141731157285696
